$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orlando")

# Copy the formatting (styles/borders) from the existing table (A1:C11) down to the new block (A13:C23)
$src = $ws.Range("A1:C11")
$src.Copy()
$dst = $ws.Range("A13:C23")
$dst.PasteSpecial(-4122)  ## xlPasteFormats
$excel.CutCopyMode = 0

# Header row
$ws.Range("A13").Value = "Year"
$ws.Range("B13").Value = "Start Day"
$ws.Range("C13").Value = "End Day"

# First data row (non-shared formula, like the original table's row 2)
$ws.Range("A14").Value = 2014
$ws.Range("B14").Formula = '="March"&" "&ROUNDUP(14-MOD((1+A14*5/4),7),0)'
$ws.Range("C14").Formula = '="November"&" "&ROUNDUP(7-MOD((1+A14*5/4),7),0)'

# Remaining years
for ($i = 1; $i -le 9; $i++) {
    $row = 14 + $i
    $year = 2014 + $i
    $ws.Cells.Item($row, 1).Value = $year
}

# Shared formulas across the rest of the block (rows 15-23), mirroring the original table
$ws.Range("B15:B23").Formula = '="March"&" "&ROUNDUP(14-MOD((1+A15*5/4),7),0)'
$ws.Range("C15:C23").Formula = '="November"&" "&ROUNDUP(7-MOD((1+A15*5/4),7),0)'

# Restore selection as in the target workbook
$ws.Range("E13").Select()
